# edit.ps1 - applies the optimization_results.xlsx changes described in the commit
# "code can now recognize no bid bids - transitions not working correctly"
#
# Summary of changes:
#  - Results sheet: Bid ID 1 has no supplier bids -> zeroed out / marked "No Bid";
#    remaining bid rows (previously split across rows 2-12) are recomputed and
#    shifted up by one row, and the now-unused last row is removed
#    (dimension shrinks from A1:O12 to A1:O11).
#  - Feasibility Notes sheet: appended a note about the excluded Bid ID.
#  - LP Model sheet: regenerated LP text reflecting the new NonBid_*_1 constraints
#    (replacing the old BidExclusion_0_1_* constraints) and simplified Transition_1_*
#    constraints now that Demand_1 is 0.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1) Results sheet: write the final values for every data cell in rows 2-11
# ---------------------------------------------------------------------------
$ws = $wb.Worksheets.Item("Results")

    # Row 2
    $ws.Range("A2").Value = 1
    $ws.Range("B2").Value = 'A'
    $ws.Range("C2").Value = 'Facility 1'
    $ws.Range("D2").Value = 'A'
    $ws.Range("E2").Value = 100
    $ws.Range("F2").Value = 0
    $ws.Range("G2").Value = 'No Bid'
    $ws.Range("H2").Value = 0
    $ws.Range("I2").Value = '''0%'
    $ws.Range("J2").Value = 0
    $ws.Range("K2").Value = 0
    $ws.Range("L2").Value = 0
    $ws.Range("M2").Value = 0
    $ws.Range("N2").Value = '''0%'
    $ws.Range("O2").Value = 0

    # Row 3
    $ws.Range("A3").Value = 2
    $ws.Range("B3").Value = 'A'
    $ws.Range("C3").Value = 'Facility 1'
    $ws.Range("D3").Value = 'B'
    $ws.Range("E3").Value = 156
    $ws.Range("F3").Value = 1404000
    $ws.Range("G3").Value = 'B'
    $ws.Range("H3").Value = 70
    $ws.Range("I3").Value = '''3%'
    $ws.Range("J3").Value = 67.89999999999999
    $ws.Range("K3").Value = 611099.9999999999
    $ws.Range("L3").Value = 9000
    $ws.Range("M3").Value = 792900.0000000001
    $ws.Range("N3").Value = '''5%'
    $ws.Range("O3").Value = 30555

    # Row 4
    $ws.Range("A4").Value = 3
    $ws.Range("B4").Value = 'A'
    $ws.Range("C4").Value = 'Facility 4'
    $ws.Range("D4").Value = 'C'
    $ws.Range("E4").Value = 423
    $ws.Range("F4").Value = 253800
    $ws.Range("G4").Value = 'C'
    $ws.Range("H4").Value = 60
    $ws.Range("I4").Value = '''4%'
    $ws.Range("J4").Value = 57.59999999999999
    $ws.Range("K4").Value = 34560
    $ws.Range("L4").Value = 600
    $ws.Range("M4").Value = 219240
    $ws.Range("N4").Value = '''7%'
    $ws.Range("O4").Value = 2419.2

    # Row 5
    $ws.Range("A5").Value = 4
    $ws.Range("B5").Value = 'A'
    $ws.Range("C5").Value = 'Facility 4'
    $ws.Range("D5").Value = 'C'
    $ws.Range("E5").Value = 453
    $ws.Range("F5").Value = 2568510
    $ws.Range("G5").Value = 'C'
    $ws.Range("H5").Value = 24
    $ws.Range("I5").Value = '''4%'
    $ws.Range("J5").Value = 23.04
    $ws.Range("K5").Value = 130636.8
    $ws.Range("L5").Value = 5670
    $ws.Range("M5").Value = 2437873.2
    $ws.Range("N5").Value = '''7%'
    $ws.Range("O5").Value = 9144.576000000001

    # Row 6
    $ws.Range("A6").Value = 5
    $ws.Range("B6").Value = 'A'
    $ws.Range("C6").Value = 'Facility 5'
    $ws.Range("D6").Value = 'C'
    $ws.Range("E6").Value = 342
    $ws.Range("F6").Value = 15390
    $ws.Range("G6").Value = 'B'
    $ws.Range("H6").Value = 34
    $ws.Range("I6").Value = '''3%'
    $ws.Range("J6").Value = 32.98
    $ws.Range("K6").Value = 1484.1
    $ws.Range("L6").Value = 45
    $ws.Range("M6").Value = 13905.9
    $ws.Range("N6").Value = '''5%'
    $ws.Range("O6").Value = 74.205

    # Row 7
    $ws.Range("A7").Value = 6
    $ws.Range("B7").Value = 'A'
    $ws.Range("C7").Value = 'Facility 9'
    $ws.Range("D7").Value = 'C'
    $ws.Range("E7").Value = 653
    $ws.Range("F7").Value = 158026
    $ws.Range("G7").Value = 'B'
    $ws.Range("H7").Value = 24
    $ws.Range("I7").Value = '''3%'
    $ws.Range("J7").Value = 23.28
    $ws.Range("K7").Value = 5633.76
    $ws.Range("L7").Value = 242
    $ws.Range("M7").Value = 152392.24
    $ws.Range("N7").Value = '''5%'
    $ws.Range("O7").Value = 281.688

    # Row 8
    $ws.Range("A8").Value = 7
    $ws.Range("B8").Value = 'A'
    $ws.Range("C8").Value = 'Facility 9'
    $ws.Range("D8").Value = 'C'
    $ws.Range("E8").Value = 432
    $ws.Range("F8").Value = 286848
    $ws.Range("G8").Value = 'A'
    $ws.Range("H8").Value = 23
    $ws.Range("I8").Value = '''0%'
    $ws.Range("J8").Value = 23
    $ws.Range("K8").Value = 15272
    $ws.Range("L8").Value = 664
    $ws.Range("M8").Value = 271576
    $ws.Range("N8").Value = '''0%'
    $ws.Range("O8").Value = 0

    # Row 9
    $ws.Range("A9").Value = 8
    $ws.Range("B9").Value = 'A'
    $ws.Range("C9").Value = 'Facility 9'
    $ws.Range("D9").Value = 'C'
    $ws.Range("E9").Value = 456
    $ws.Range("F9").Value = 10944
    $ws.Range("G9").Value = 'B'
    $ws.Range("H9").Value = 13
    $ws.Range("I9").Value = '''3%'
    $ws.Range("J9").Value = 12.61
    $ws.Range("K9").Value = 302.64
    $ws.Range("L9").Value = 24
    $ws.Range("M9").Value = 10641.36
    $ws.Range("N9").Value = '''5%'
    $ws.Range("O9").Value = 15.132

    # Row 10
    $ws.Range("A10").Value = 9
    $ws.Range("B10").Value = 'A'
    $ws.Range("C10").Value = 'Facility 10'
    $ws.Range("D10").Value = 'C'
    $ws.Range("E10").Value = 234
    $ws.Range("F10").Value = 54288
    $ws.Range("G10").Value = 'C'
    $ws.Range("H10").Value = 32
    $ws.Range("I10").Value = '''4%'
    $ws.Range("J10").Value = 30.72
    $ws.Range("K10").Value = 7127.04
    $ws.Range("L10").Value = 232
    $ws.Range("M10").Value = 47160.96
    $ws.Range("N10").Value = '''7%'
    $ws.Range("O10").Value = 498.8928

    # Row 11
    $ws.Range("A11").Value = 10
    $ws.Range("B11").Value = 'A'
    $ws.Range("C11").Value = 'Facility 10'
    $ws.Range("D11").Value = 'C'
    $ws.Range("E11").Value = 231
    $ws.Range("F11").Value = 3003
    $ws.Range("G11").Value = 'B'
    $ws.Range("H11").Value = 13
    $ws.Range("I11").Value = '''3%'
    $ws.Range("J11").Value = 12.61
    $ws.Range("K11").Value = 163.93
    $ws.Range("L11").Value = 13
    $ws.Range("M11").Value = 2839.07
    $ws.Range("N11").Value = '''5%'
    $ws.Range("O11").Value = 8.1965

    # Remove the now-obsolete last row (old row 12); shrinks used range to A1:O11
    $ws.Rows.Item(12).Delete()

# ---------------------------------------------------------------------------
# 2) Feasibility Notes sheet: append note about the excluded Bid ID
# ---------------------------------------------------------------------------
$wsNotes = $wb.Worksheets.Item("Feasibility Notes")
$wsNotes.Range("A2").Value = "Model is optimal." + "`n" + "Note: The following Bid ID(s) were excluded because no supplier bids were found: 1."

# ---------------------------------------------------------------------------
# 3) LP Model sheet: replace the generated LP text with the regenerated model
# ---------------------------------------------------------------------------
$wsLP = $wb.Worksheets.Item("LP Model")
$lpText = @'
\* Sourcing_with_MultiTier_Rebates_Discounts *\
Minimize
OBJ: S_A + S_B + S_C - rebate_A - rebate_B - rebate_C
Subject To
BaseSpend_A: S0_A - 64 x_A_10 - 70 x_A_2 - 55 x_A_3 - 23 x_A_4 - 54 x_A_5
 - 42 x_A_6 - 23 x_A_7 - 75 x_A_8 - 97 x_A_9 = 0
BaseSpend_B: S0_B - 13 x_B_10 - 70 x_B_2 - 65 x_B_3 - 75 x_B_4 - 34 x_B_5
 - 24 x_B_6 - 85 x_B_7 - 13 x_B_8 - 56 x_B_9 = 0
BaseSpend_C: S0_C - 15 x_C_10 - 75 x_C_2 - 60 x_C_3 - 24 x_C_4 - 44 x_C_6
 - 42 x_C_7 - 24 x_C_8 - 32 x_C_9 = 0
Capacity_B_Bid_ID_1: x_B_1 <= 100000000
Capacity_B_Bid_ID_10: x_B_10 <= 100000000
Capacity_B_Bid_ID_2: x_B_2 <= 100000000
Capacity_B_Bid_ID_3: x_B_3 <= 100000000
Capacity_B_Bid_ID_4: x_B_4 <= 100000000
Capacity_B_Bid_ID_5: x_B_5 <= 100000000
Capacity_B_Bid_ID_6: x_B_6 <= 100000000
Capacity_B_Bid_ID_7: x_B_7 <= 100000000
Capacity_B_Bid_ID_8: x_B_8 <= 100000000
Capacity_B_Bid_ID_9: x_B_9 <= 100000000
Capacity_C_Bid_ID_1: x_C_1 <= 100000000
Capacity_C_Bid_ID_10: x_C_10 <= 100000000
Capacity_C_Bid_ID_2: x_C_2 <= 100000000
Capacity_C_Bid_ID_3: x_C_3 <= 100000000
Capacity_C_Bid_ID_4: x_C_4 <= 100000000
Capacity_C_Bid_ID_5: x_C_5 <= 100000000
Capacity_C_Bid_ID_6: x_C_6 <= 100000000
Capacity_C_Bid_ID_7: x_C_7 <= 100000000
Capacity_C_Bid_ID_8: x_C_8 <= 100000000
Capacity_C_Bid_ID_9: x_C_9 <= 100000000
Demand_1: x_A_1 + x_B_1 + x_C_1 = 0
Demand_10: x_A_10 + x_B_10 + x_C_10 = 13
Demand_2: x_A_2 + x_B_2 + x_C_2 = 9000
Demand_3: x_A_3 + x_B_3 + x_C_3 = 600
Demand_4: x_A_4 + x_B_4 + x_C_4 = 5670
Demand_5: x_A_5 + x_B_5 + x_C_5 = 45
Demand_6: x_A_6 + x_B_6 + x_C_6 = 242
Demand_7: x_A_7 + x_B_7 + x_C_7 = 664
Demand_8: x_A_8 + x_B_8 + x_C_8 = 24
Demand_9: x_A_9 + x_B_9 + x_C_9 = 232
DiscountTierLower_A_0: d_A - 19400000000 z_discount_A_0 >= -19400000000
DiscountTierLower_A_1: - 0.01 S0_A + d_A - 19400000000 z_discount_A_1
 >= -19400000000
DiscountTierLower_B_0: d_B - 97000000000 z_discount_B_0 >= -97000000000
DiscountTierLower_B_1: - 0.03 S0_B + d_B - 97000000000 z_discount_B_1
 >= -97000000000
DiscountTierLower_C_0: d_C - 97000000000 z_discount_C_0 >= -97000000000
DiscountTierLower_C_1: - 0.04 S0_C + d_C - 97000000000 z_discount_C_1
 >= -97000000000
DiscountTierMax_A_0: 19400000000 z_discount_A_0 <= 19400001000
DiscountTierMax_B_0: 97000000000 z_discount_B_0 <= 97000000500
DiscountTierMax_C_0: 97000000000 z_discount_C_0 <= 97000000500
_dummy: __dummy = 0
DiscountTierMin_A_0: __dummy >= 0
DiscountTierMin_A_1: x_A_1 + x_A_10 + x_A_3 + x_A_4 + x_A_8 + x_A_9
 - 1000 z_discount_A_1 >= 0
DiscountTierMin_B_0: __dummy >= 0
DiscountTierMin_B_1: x_B_2 + x_B_5 + x_B_6 + x_B_7 - 500 z_discount_B_1 >= 0
DiscountTierMin_C_0: __dummy >= 0
DiscountTierMin_C_1: x_C_1 + x_C_10 + x_C_3 + x_C_4 + x_C_8 + x_C_9
 - 500 z_discount_C_1 >= 0
DiscountTierSelect_A: z_discount_A_0 + z_discount_A_1 = 1
DiscountTierSelect_B: z_discount_B_0 + z_discount_B_1 = 1
DiscountTierSelect_C: z_discount_C_0 + z_discount_C_1 = 1
DiscountTierUpper_A_0: d_A + 19400000000 z_discount_A_0 <= 19400000000
DiscountTierUpper_A_1: - 0.01 S0_A + d_A + 19400000000 z_discount_A_1
 <= 19400000000
DiscountTierUpper_B_0: d_B + 97000000000 z_discount_B_0 <= 97000000000
DiscountTierUpper_B_1: - 0.03 S0_B + d_B + 97000000000 z_discount_B_1
 <= 97000000000
DiscountTierUpper_C_0: d_C + 97000000000 z_discount_C_0 <= 97000000000
DiscountTierUpper_C_1: - 0.04 S0_C + d_C + 97000000000 z_discount_C_1
 <= 97000000000
EffectiveSpend_A: - S0_A + S_A + d_A = 0
EffectiveSpend_B: - S0_B + S_B + d_B = 0
EffectiveSpend_C: - S0_C + S_C + d_C = 0
NonBid_A_1: x_A_1 = 0
NonBid_B_1: x_B_1 = 0
NonBid_C_1: x_C_1 = 0
NonBid_C_5: x_C_5 = 0
RebateTierLower_A_0: rebate_A - 19400000000 y_rebate_A_0 >= -19400000000
RebateTierLower_A_1: - 0.1 S_A + rebate_A - 19400000000 y_rebate_A_1
 >= -19400000000
RebateTierLower_B_0: rebate_B - 97000000000 y_rebate_B_0 >= -97000000000
RebateTierLower_B_1: - 0.05 S_B + rebate_B - 97000000000 y_rebate_B_1
 >= -97000000000
RebateTierLower_C_0: rebate_C - 97000000000 y_rebate_C_0 >= -97000000000
RebateTierLower_C_1: - 0.07 S_C + rebate_C - 97000000000 y_rebate_C_1
 >= -97000000000
RebateTierMax_A_0: 19400000000 y_rebate_A_0 <= 19400000500
RebateTierMax_B_0: 97000000000 y_rebate_B_0 <= 97000000500
RebateTierMax_C_0: 97000000000 y_rebate_C_0 <= 97000000700
RebateTierMin_A_0: __dummy >= 0
RebateTierMin_A_1: - 500 y_rebate_A_1 >= 0
RebateTierMin_B_0: __dummy >= 0
RebateTierMin_B_1: x_B_2 + x_B_5 + x_B_6 + x_B_7 - 500 y_rebate_B_1 >= 0
RebateTierMin_C_0: __dummy >= 0
RebateTierMin_C_1: x_C_1 + x_C_10 + x_C_3 + x_C_4 + x_C_8 + x_C_9
 - 700 y_rebate_C_1 >= 0
RebateTierSelect_A: y_rebate_A_0 + y_rebate_A_1 = 1
RebateTierSelect_B: y_rebate_B_0 + y_rebate_B_1 = 1
RebateTierSelect_C: y_rebate_C_0 + y_rebate_C_1 = 1
RebateTierUpper_A_0: rebate_A + 19400000000 y_rebate_A_0 <= 19400000000
RebateTierUpper_A_1: - 0.1 S_A + rebate_A + 19400000000 y_rebate_A_1
 <= 19400000000
RebateTierUpper_B_0: rebate_B + 97000000000 y_rebate_B_0 <= 97000000000
RebateTierUpper_B_1: - 0.05 S_B + rebate_B + 97000000000 y_rebate_B_1
 <= 97000000000
RebateTierUpper_C_0: rebate_C + 97000000000 y_rebate_C_0 <= 97000000000
RebateTierUpper_C_1: - 0.07 S_C + rebate_C + 97000000000 y_rebate_C_1
 <= 97000000000
Transition_10_A: - 13 T_10_A + x_A_10 <= 0
Transition_10_B: - 13 T_10_B + x_B_10 <= 0
Transition_1_B: x_B_1 <= 0
Transition_1_C: x_C_1 <= 0
Transition_2_A: - 9000 T_2_A + x_A_2 <= 0
Transition_2_C: - 9000 T_2_C + x_C_2 <= 0
Transition_3_A: - 600 T_3_A + x_A_3 <= 0
Transition_3_B: - 600 T_3_B + x_B_3 <= 0
Transition_4_A: - 5670 T_4_A + x_A_4 <= 0
Transition_4_B: - 5670 T_4_B + x_B_4 <= 0
Transition_5_A: - 45 T_5_A + x_A_5 <= 0
Transition_5_B: - 45 T_5_B + x_B_5 <= 0
Transition_6_A: - 242 T_6_A + x_A_6 <= 0
Transition_6_B: - 242 T_6_B + x_B_6 <= 0
Transition_7_A: - 664 T_7_A + x_A_7 <= 0
Transition_7_B: - 664 T_7_B + x_B_7 <= 0
Transition_8_A: - 24 T_8_A + x_A_8 <= 0
Transition_8_B: - 24 T_8_B + x_B_8 <= 0
Transition_9_A: - 232 T_9_A + x_A_9 <= 0
Transition_9_B: - 232 T_9_B + x_B_9 <= 0
Volume_A: V_A - x_A_1 - x_A_10 - x_A_2 - x_A_3 - x_A_4 - x_A_5 - x_A_6 - x_A_7
 - x_A_8 - x_A_9 = 0
Volume_B: V_B - x_B_1 - x_B_10 - x_B_2 - x_B_3 - x_B_4 - x_B_5 - x_B_6 - x_B_7
 - x_B_8 - x_B_9 = 0
Volume_C: V_C - x_C_1 - x_C_10 - x_C_2 - x_C_3 - x_C_4 - x_C_5 - x_C_6 - x_C_7
 - x_C_8 - x_C_9 = 0
Binaries
T_10_A
T_10_B
T_2_A
T_2_C
T_3_A
T_3_B
T_4_A
T_4_B
T_5_A
T_5_B
T_6_A
T_6_B
T_7_A
T_7_B
T_8_A
T_8_B
T_9_A
T_9_B
y_rebate_A_0
y_rebate_A_1
y_rebate_B_0
y_rebate_B_1
y_rebate_C_0
y_rebate_C_1
z_discount_A_0
z_discount_A_1
z_discount_B_0
z_discount_B_1
z_discount_C_0
z_discount_C_1
End
'@ + "`n"
$wsLP.Range("A2").Value = $lpText
